# Generate Report for Handoff
# Adds two new localization entries (f1dee5f1-... and f6ffeb92-...) as rows 6/7
# on all three sheets: Overview, zh-cn, de-de. Mirrors the existing
# "Ready for handoff" rows (e.g. dca1bbfd-...) already present as row 5.

$wb = $excel.ActiveWorkbook

$repo = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob"
$commit1 = "79b6e5f1c6a4d9b2e3f4a5b6c7d8e9f0a1b2c3d4"
$commit2 = "8a7b6c5d4e3f2a1b0c9d8e7f6a5b4c3d2e1f0a9b"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns A..G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "f1dee5f1-ab06-4332-9fce-140a7e42c38d.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-09-04 10:45:49"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "$repo/$commit1/e2e/f1dee5f1-ab06-4332-9fce-140a7e42c38d.md", "", "", "e2e\f1dee5f1-ab06-4332-9fce-140a7e42c38d.md")

$wsOverview.Range("A7").Value = "f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-04 10:45:49"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "$repo/$commit2/e2e/f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md", "", "", "e2e\f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md")

$wsOverview.Range("G6:G7").NumberFormat = $wsOverview.Range("G5").NumberFormat

$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G7"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - columns A..P
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "f1dee5f1-ab06-4332-9fce-140a7e42c38d.b12a6462b2526bd4568e28f44e89fc38ef7a332b.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-09-04 10:45:44"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "$repo/$commit1/e2e/f1dee5f1-ab06-4332-9fce-140a7e42c38d.md", "", "", "f1dee5f1-ab06-4332-9fce-140a7e42c38d.md")

$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "f6ffeb92-9edf-4d59-b02d-ed3977ac2746.7860a6873ce8b46d3ac50fe5fe7fa4acfb37a591.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-04 10:45:44"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "$repo/$commit2/e2e/f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md", "", "", "f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md")

$wsZhCn.Range("H6:H7").NumberFormat = $wsZhCn.Range("H5").NumberFormat
$wsZhCn.Range("K6:K7").NumberFormat = $wsZhCn.Range("K5").NumberFormat

$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P7"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - columns A..P
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "f1dee5f1-ab06-4332-9fce-140a7e42c38d.b12a6462b2526bd4568e28f44e89fc38ef7a332b.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-09-04 10:45:49"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "$repo/$commit1/e2e/f1dee5f1-ab06-4332-9fce-140a7e42c38d.md", "", "", "f1dee5f1-ab06-4332-9fce-140a7e42c38d.md")

$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "f6ffeb92-9edf-4d59-b02d-ed3977ac2746.7860a6873ce8b46d3ac50fe5fe7fa4acfb37a591.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-04 10:45:49"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "$repo/$commit2/e2e/f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md", "", "", "f6ffeb92-9edf-4d59-b02d-ed3977ac2746.md")

$wsDeDe.Range("H6:H7").NumberFormat = $wsDeDe.Range("H5").NumberFormat
$wsDeDe.Range("K6:K7").NumberFormat = $wsDeDe.Range("K5").NumberFormat

$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P7"))
